$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): A1:C1 already carry the bold/bordered/centered
#     header style from the template; fill in their text and extend the
#     same header style across the new D1:I1 columns.
$ws.Range("A1").Value = "Job_Id"
$ws.Range("B1").Value = "Job_Title"
$ws.Range("C1").Value = "Job_Description"
$ws.Range("D1").Value = "Total_Years_Min_Exp"
$ws.Range("E1").Value = "Total_Years_Max_Exp"
$ws.Range("F1").Value = "LinkedIn_Poster"
$ws.Range("G1").Value = "LinkedIn_Posted"
$ws.Range("H1").Value = "Resume_received"
$ws.Range("I1").Value = "Resume_downloaded"

# Copy the header formatting from A1 onto the newly added header cells so
# they reuse the very same cell style (bold font, thin border, centered).
$ws.Range("A1").Copy()
$ws.Range("D1:I1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data row (row 2)
$ws.Range("A2").ClearFormats()
$ws.Range("A2").Value = "JD_001"
$ws.Range("B2").Value = "Senior Engineer"
$ws.Range("C2").Value = "We are seeking a Software Engineer to build and maintain high-quality software solutions.`nWork with global teams to drive innovation and deliver scalable applications.`nJoin Akkodis and be part of a tech-driven, collaborative environment."
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 2

# Recompute row height so the embedded line breaks in C2 don't leave a
# stale custom row height behind.
$ws.Rows(2).AutoFit()
